# fix: correção na exportação de planilha da Estante Virtual para o Bling
#
# The "Estante Virtual" export template gains a new "Nº Checkout" column
# (inserted as the 3rd column, right after "Status") and the existing
# "E-mail do Comprador" column (previously right after "Comprador") is
# relocated to become the very last column of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the "E-mail do Comprador" column from its current position
#    (column L). Everything to its right (M:AG) shifts left by one.
$ws.Range("L1").EntireColumn.Delete()

# 2) Insert a new, blank column at C (shifting the old "Pedido.." block,
#    now C:AG, right by one) and give its header the new field name.
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").Value = "Nº Checkout"

# 3) Re-add "E-mail do Comprador" as the new last column (AH).
$ws.Range("AH1").Value = "E-mail do Comprador"
